# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly scraped counts (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rows keyed by F column) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 320
$ws1.Range("F4").Value = 245
$ws1.Range("F5").Value = 2904
$ws1.Range("F6").Value = 1991
$ws1.Range("F7").Value = 384
$ws1.Range("F9").Value = 1065
$ws1.Range("F10").Value = 201
$ws1.Range("F11").Value = 421
$ws1.Range("F12").Value = 51

# --- Sheet "全部类型" (rows keyed by F column, one extra row vs 展览) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 320
$ws4.Range("F4").Value = 245
$ws4.Range("F5").Value = 2904
$ws4.Range("F6").Value = 1991
$ws4.Range("F7").Value = 384
$ws4.Range("F10").Value = 1065
$ws4.Range("F11").Value = 201
$ws4.Range("F12").Value = 422
$ws4.Range("F13").Value = 51
